$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# are preserved exactly as text (not coerced into floating point numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.873.68"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "1.642.70"
$ws.Range("E3").Value = "  +0.77%  "

$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.42%  "

$ws.Range("D5").Value = "215.83"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").Value = "0.5057"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").Value = "0.2589"
$ws.Range("E8").Value = "  +0.61%  "

$ws.Range("D9").Value = "0.06436"
$ws.Range("E9").Value = "  +1.79%  "

$ws.Range("D10").Value = "20.53"
$ws.Range("E10").Value = "  +5.43%  "

$ws.Range("D11").Value = "0.07809"
$ws.Range("E11").Value = "  +0.78%  "

$ws.Range("D12").Value = "4.280"
$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("D13").Value = "1.866.31"
$ws.Range("E13").Value = "  +0.60%  "

$ws.Range("D14").Value = "1.625.48"
$ws.Range("E14").Value = "  -0.73%  "

$ws.Range("D15").Value = "0.5617"
$ws.Range("E15").Value = "  +2.09%  "

$ws.Range("D16").Value = "0.0₅7696"
$ws.Range("E16").Value = "  +0.60%  "

$ws.Range("D17").Value = "63.26"
$ws.Range("E17").Value = "  -0.63%  "

$ws.Range("D18").Value = "25.870.71"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").Value = "193.75"
$ws.Range("E20").Value = "  -0.29%  "

$ws.Range("D21").Value = "4.367"
$ws.Range("E21").Value = "  -1.01%  "

$ws.Range("D22").Value = "9.943"
$ws.Range("E22").Value = "  +0.79%  "

$ws.Range("D23").Value = "6.125"
$ws.Range("E23").Value = "  +1.58%  "

$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").Value = "1.798"
$ws.Range("E25").Value = "  -6.01%  "

$ws.Range("D26").Value = "140.47"
$ws.Range("E26").Value = "  -1.07%  "

$ws.Range("D27").Value = "0.1238"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").Value = "6.823"
$ws.Range("E28").Value = "  +0.84%  "

$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").Value = "1.245"
$ws.Range("E30").Value = "  +0.64%  "

$ws.Range("D31").Value = "0.04964"
$ws.Range("E31").Value = "  +1.65%  "

$ws.Range("D32").Value = "3.296"
$ws.Range("E32").Value = "  +1.59%  "

$ws.Range("D33").Value = "3.237"
$ws.Range("E33").Value = "  +1.46%  "

$ws.Range("D34").Value = "1.579"
$ws.Range("E34").Value = "  +2.28%  "

$ws.Range("D35").Value = "2.385"
$ws.Range("E35").Value = "  +0.72%  "

$ws.Range("D36").Value = "0.9054"
$ws.Range("E36").Value = "  +1.15%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "2.573"

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.5575"
$ws.Range("E38").Value = "  +0.90%  "

$ws.Range("D39").Value = "1.128.81"
$ws.Range("E39").Value = "  +0.67%  "

$ws.Range("D40").Value = "0.01565"
$ws.Range("E40").Value = "  +1.06%  "

$ws.Range("D41").Value = "0.9969"
$ws.Range("E41").Value = "  -0.37%  "

$ws.Range("D42").Value = "5.521"
$ws.Range("E42").Value = "  -1.09%  "

$ws.Range("D43").Value = "0.8026"
$ws.Range("E43").Value = "  +0.96%  "

$ws.Range("D44").Value = "98.66"
$ws.Range("E44").Value = "  +1.45%  "

$ws.Range("D45").Value = "1.778.26"
$ws.Range("E45").Value = "  +0.78%  "

$ws.Range("D46").Value = "0.0₈111"
$ws.Range("E46").Value = "  -6.45%  "

$ws.Range("D47").Value = "55.70"
$ws.Range("E47").Value = "  +1.84%  "

$ws.Range("D48").Value = "0.4285"
$ws.Range("E48").Value = "  -3.54%  "

$ws.Range("D49").Value = "7.736"
$ws.Range("E49").Value = "  +1.90%  "

$ws.Range("D50").Value = "0.05049"
$ws.Range("E50").Value = "  -1.71%  "

$ws.Range("D51").Value = "0.9956"
$ws.Range("E51").Value = "  -0.92%  "
